$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 7.892066666666667
$ws.Range("H2").Value = 23.6762
$ws.Range("I2").Value = 0.1739002798877711
$ws.Range("J2").Value = 0.1739002798877711
$ws.Range("M2").Value = 0.5001966666666666
$ws.Range("N2").Value = 1.50059
$ws.Range("O2").Value = 0.03894027965151046
$ws.Range("P2").Value = 0.03894027965151046
$ws.Range("Q2").Value = 3.947585439777778
$ws.Range("R2").Value = 35.528268958
$ws.Range("S2").Value = 0.006771725530305747
$ws.Range("T2").Value = 0.006771725530305747

$ws.Range("G3").Value = 7.892066666666667
$ws.Range("H3").Value = 23.6762
$ws.Range("I3").Value = 0.1739002798877711
$ws.Range("J3").Value = 0.1739002798877711
$ws.Range("O3").Value = 0.7732779360092192
$ws.Range("P3").Value = 0.7732779360092191
$ws.Range("Q3").Value = 78.39134049395555
$ws.Range("R3").Value = 705.5220644456
$ws.Range("S3").Value = 0.1344732495030412
$ws.Range("T3").Value = 0.1344732495030411

$ws.Range("G4").Value = 7.892066666666667
$ws.Range("H4").Value = 23.6762
$ws.Range("I4").Value = 0.1739002798877711
$ws.Range("J4").Value = 0.1739002798877711
$ws.Range("M4").Value = 2.334238666666666
$ws.Range("N4").Value = 7.002715999999999
$ws.Range("O4").Value = 0.1817203362411497
$ws.Range("P4").Value = 0.1817203362411496
$ws.Range("Q4").Value = 18.42196717324444
$ws.Range("R4").Value = 165.7977045592
$ws.Range("S4").Value = 0.03160121733363579
$ws.Range("T4").Value = 0.03160121733363579

$ws.Range("G5").Value = 7.892066666666667
$ws.Range("H5").Value = 23.6762
$ws.Range("I5").Value = 0.1739002798877711
$ws.Range("J5").Value = 0.1739002798877711
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.07786066666666666
$ws.Range("N5").Value = 0.233582
$ws.Range("O5").Value = 0.006061448098120818
$ws.Range("P5").Value = 0.006061448098120817
$ws.Range("Q5").Value = 0.6144815720444444
$ws.Range("R5").Value = 5.5303341484
$ws.Range("S5").Value = 0.001054087520788408
$ws.Range("T5").Value = 0.001054087520788408

$ws.Range("I6").Value = 0.3815924715300191
$ws.Range("J6").Value = 0.3815924715300191
$ws.Range("M6").Value = 0.5001966666666666
$ws.Range("N6").Value = 1.50059
$ws.Range("O6").Value = 0.03894027965151046
$ws.Range("P6").Value = 0.03894027965151046
$ws.Range("Q6").Value = 8.662256814726666
$ws.Range("R6").Value = 77.96031133254
$ws.Range("S6").Value = 0.01485931755428999
$ws.Range("T6").Value = 0.01485931755428999

$ws.Range("I7").Value = 0.3815924715300191
$ws.Range("J7").Value = 0.3815924715300191
$ws.Range("O7").Value = 0.7732779360092192
$ws.Range("P7").Value = 0.7732779360092191
$ws.Range("S7").Value = 0.2950770387813899
$ws.Range("T7").Value = 0.2950770387813899

$ws.Range("I8").Value = 0.3815924715300191
$ws.Range("J8").Value = 0.3815924715300191
$ws.Range("M8").Value = 2.334238666666666
$ws.Range("N8").Value = 7.002715999999999
$ws.Range("O8").Value = 0.1817203362411497
$ws.Range("P8").Value = 0.1817203362411496
$ws.Range("Q8").Value = 40.42364962621066
$ws.Range("R8").Value = 363.812846635896
$ws.Range("S8").Value = 0.06934311223352639
$ws.Range("T8").Value = 0.06934311223352639

$ws.Range("I9").Value = 0.3815924715300191
$ws.Range("J9").Value = 0.3815924715300191
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.07786066666666666
$ws.Range("N9").Value = 0.233582
$ws.Range("O9").Value = 0.006061448098120818
$ws.Range("P9").Value = 0.006061448098120817
$ws.Range("Q9").Value = 1.348367822854667
$ws.Range("R9").Value = 12.135310405692
$ws.Range("S9").Value = 0.002313002960812856
$ws.Range("T9").Value = 0.002313002960812856

$ws.Range("G10").Value = 7.716272666666666
$ws.Range("H10").Value = 23.148818
$ws.Range("I10").Value = 0.1700266904854272
$ws.Range("J10").Value = 0.1700266904854272
$ws.Range("M10").Value = 0.5001966666666666
$ws.Range("N10").Value = 1.50059
$ws.Range("O10").Value = 0.03894027965151046
$ws.Range("P10").Value = 0.03894027965151046
$ws.Range("Q10").Value = 3.859653866957777
$ws.Range("R10").Value = 34.73688480262
$ws.Range("S10").Value = 0.006620886875723349
$ws.Range("T10").Value = 0.00662088687572335

$ws.Range("G11").Value = 7.716272666666666
$ws.Range("H11").Value = 23.148818
$ws.Range("I11").Value = 0.1700266904854272
$ws.Range("J11").Value = 0.1700266904854272
$ws.Range("O11").Value = 0.7732779360092192
$ws.Range("P11").Value = 0.7732779360092191
$ws.Range("Q11").Value = 76.64519111473155
$ws.Range("R11").Value = 689.8067200325839
$ws.Range("S11").Value = 0.1314778882850495
$ws.Range("T11").Value = 0.1314778882850495

$ws.Range("G12").Value = 7.716272666666666
$ws.Range("H12").Value = 23.148818
$ws.Range("I12").Value = 0.1700266904854272
$ws.Range("J12").Value = 0.1700266904854272
$ws.Range("M12").Value = 2.334238666666666
$ws.Range("N12").Value = 7.002715999999999
$ws.Range("O12").Value = 0.1817203362411497
$ws.Range("P12").Value = 0.1817203362411496
$ws.Range("Q12").Value = 18.01162202107644
$ws.Range("R12").Value = 162.104598189688
$ws.Range("S12").Value = 0.03089730736498171
$ws.Range("T12").Value = 0.03089730736498171

$ws.Range("G13").Value = 7.716272666666666
$ws.Range("H13").Value = 23.148818
$ws.Range("I13").Value = 0.1700266904854272
$ws.Range("J13").Value = 0.1700266904854272
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.07786066666666666
$ws.Range("N13").Value = 0.233582
$ws.Range("O13").Value = 0.006061448098120818
$ws.Range("P13").Value = 0.006061448098120817
$ws.Range("Q13").Value = 0.6007941340084444
$ws.Range("R13").Value = 5.407147206075999
$ws.Range("S13").Value = 0.00103060795967267
$ws.Range("T13").Value = 0.00103060795967267

$ws.Range("G14").Value = 12.45667266666667
$ws.Range("H14").Value = 37.370018
$ws.Range("I14").Value = 0.2744805580967825
$ws.Range("J14").Value = 0.2744805580967826
$ws.Range("M14").Value = 0.5001966666666666
$ws.Range("N14").Value = 1.50059
$ws.Range("O14").Value = 0.03894027965151046
$ws.Range("P14").Value = 0.03894027965151046
$ws.Range("Q14").Value = 6.230786145624444
$ws.Range("R14").Value = 56.07707531062
$ws.Range("S14").Value = 0.01068834969119138
$ws.Range("T14").Value = 0.01068834969119138

$ws.Range("G15").Value = 12.45667266666667
$ws.Range("H15").Value = 37.370018
$ws.Range("I15").Value = 0.2744805580967825
$ws.Range("J15").Value = 0.2744805580967826
$ws.Range("O15").Value = 0.7732779360092192
$ws.Range("P15").Value = 0.7732779360092191
$ws.Range("Q15").Value = 123.7312493264649
$ws.Range("R15").Value = 1113.581243938184
$ws.Range("S15").Value = 0.2122497594397386
$ws.Range("T15").Value = 0.2122497594397386

$ws.Range("G16").Value = 12.45667266666667
$ws.Range("H16").Value = 37.370018
$ws.Range("I16").Value = 0.2744805580967825
$ws.Range("J16").Value = 0.2744805580967826
$ws.Range("M16").Value = 2.334238666666666
$ws.Range("N16").Value = 7.002715999999999
$ws.Range("O16").Value = 0.1817203362411497
$ws.Range("P16").Value = 0.1817203362411496
$ws.Range("Q16").Value = 29.07684699654311
$ws.Range("R16").Value = 261.691622968888
$ws.Range("S16").Value = 0.04987869930900574
$ws.Range("T16").Value = 0.04987869930900574

$ws.Range("G17").Value = 12.45667266666667
$ws.Range("H17").Value = 37.370018
$ws.Range("I17").Value = 0.2744805580967825
$ws.Range("J17").Value = 0.2744805580967826
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.07786066666666666
$ws.Range("N17").Value = 0.233582
$ws.Range("O17").Value = 0.006061448098120818
$ws.Range("P17").Value = 0.006061448098120817
$ws.Range("Q17").Value = 0.9698848382751112
$ws.Range("R17").Value = 8.728963544476001
$ws.Range("S17").Value = 0.001663749656846883
$ws.Range("T17").Value = 0.001663749656846883
